# This script reorders the data rows (rows 2..75) of the single worksheet
# according to the new order produced by the upstream export/log process,
# and bumps the "Förändrad" (column C) date serial for every data row
# from 46077 to 46078. No cell content is invented: every value/formula
# already exists somewhere in the original sheet; rows are simply moved
# to their new positions (a couple of rows also pick up/lose an
# already-present "Markägare" (F) value and trailing link columns, which
# is captured automatically because we copy the *entire* source row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstDataRow = 2
$lastDataRow = 75
$lastCol = 26   # column Z
$newChanged = 46078

# Target order of "Beteckning" (column A) values for rows 2..75, top to bottom.
$targetOrder = @(
  "A 43913-2021",
  "A 695-2023",
  "A 48514-2025",
  "A 60416-2025",
  "A 24599-2022",
  "A 45595-2025",
  "A 43900-2021",
  "A 14220-2025",
  "A 1697-2023",
  "A 30067-2025",
  "A 30055-2025",
  "A 60406-2025",
  "A 56816-2021",
  "A 14214-2025",
  "A 23829-2023",
  "A 33801-2025",
  "A 3726-2022",
  "A 1691-2023",
  "A 18100-2022",
  "A 27791-2021",
  "A 3681-2022",
  "A 48667-2022",
  "A 32828-2022",
  "A 43717-2021",
  "A 33015-2021",
  "A 33016-2021",
  "A 20210-2021",
  "A 73772-2021",
  "A 10817-2021",
  "A 52035-2021",
  "A 35492-2021",
  "A 48647-2022",
  "A 14202-2025",
  "A 21285-2023",
  "A 61242-2022",
  "A 30925-2023",
  "A 15690-2023",
  "A 16142-2024",
  "A 14217-2025",
  "A 11876-2024",
  "A 30378-2024",
  "A 61320-2024",
  "A 73194-2021",
  "A 49175-2025",
  "A 45945-2025",
  "A 15905-2022",
  "A 16483-2023",
  "A 50474-2025",
  "A 46405-2022",
  "A 20649-2022",
  "A 8528-2024",
  "A 45869-2022",
  "A 3686-2022",
  "A 54468-2025",
  "A 15692-2023",
  "A 73151-2021",
  "A 9032-2023",
  "A 73155-2021",
  "A 48656-2022",
  "A 38194-2022",
  "A 11867-2024",
  "A 16149-2024",
  "A 28983-2024",
  "A 61225-2022",
  "A 25217-2024",
  "A 35197-2024",
  "A 30622-2024",
  "A 36097-2021",
  "A 16145-2024",
  "A 61514-2022",
  "A 3021-2025",
  "A 6059-2022",
  "A 37095-2021",
  "A 10815-2021"
)

# --- Step 1: snapshot every current data row (formula/value per cell) into
#     a lookup table keyed by its "Beteckning" (column A). Using .Formula
#     preserves both plain values (numbers/strings/dates) and HYPERLINK()
#     formulas alike, so no special-casing is needed per column.
$rowsByKey = @{}
for ($r = $firstDataRow; $r -le $lastDataRow; $r++) {
    $key = $ws.Cells.Item($r, 1).Formula
    $rowVals = New-Object 'object[]' $lastCol
    for ($c = 1; $c -le $lastCol; $c++) {
        $rowVals[$c - 1] = $ws.Cells.Item($r, $c).Formula
    }
    $rowsByKey[$key] = $rowVals
}

# --- Step 2: write every row back out in its new position, using the
#     snapshot captured above, and force column C to the new date serial.
for ($i = 0; $i -lt $targetOrder.Count; $i++) {
    $destRow = $firstDataRow + $i
    $key = $targetOrder[$i]
    $rowVals = $rowsByKey[$key]
    for ($c = 1; $c -le $lastCol; $c++) {
        $ws.Cells.Item($destRow, $c).Formula = $rowVals[$c - 1]
    }
    $ws.Cells.Item($destRow, 3).Formula = $newChanged
}
